# Daily attendance processing - 2025-11-25 15:29:18
# Normalize the "Recorded By" column (G) so that rows recorded jointly by
# dnasr281@gmail.com and "System" (or admin@admin.com) list the other
# party first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value = "admin@admin.com, dnasr281@gmail.com"
    }
}
